$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly "Betarraga" price sheet gets a brand-new observation inserted
# as row 130; every following record (old rows 130-177) shifts down by one
# row (to 131-178), just like a new week being prepended to the series.
$ws.Rows(130).Insert()

# Populate the newly inserted row 130 with this week's data.
$ws.Cells.Item(130, 1).Value = 5
$ws.Cells.Item(130, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(130, 3).Value = "Maule"
$ws.Cells.Item(130, 4).Value = 44468
$ws.Cells.Item(130, 5).Value = 7
$ws.Cells.Item(130, 6).Value = 100114014
$ws.Cells.Item(130, 7).Value = "Betarraga"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 3000
$ws.Cells.Item(130, 11).Value = 650
$ws.Cells.Item(130, 12).Value = 650
$ws.Cells.Item(130, 13).Value = 650
$ws.Cells.Item(130, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(130, 15).Value = "Región del Maule"
$ws.Cells.Item(130, 16).Value = 130
$ws.Cells.Item(130, 17).Value = 5
$ws.Cells.Item(130, 18).Value = "Hortaliza"
